$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44557
$ws.Range("K2").Value = 'Lapins'
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 9000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 9500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 950
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44208
$ws.Range("K3").Value = 'Lapins'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 10500
$ws.Range("O3").Value = 11000
$ws.Range("P3").Value = 10750
$ws.Range("Q3").Value = '$/bandeja 12 kilos'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 896
$ws.Range("T3").Value = 12

# Row 4
$ws.Range("D4").Value = 44175
$ws.Range("K4").Value = 'Rainier'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 26000
$ws.Range("P4").Value = 25500
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1417
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44532
$ws.Range("K5").Value = 'Brooks'
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 27000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 27500
$ws.Range("Q5").Value = '$/bandeja 12 kilos'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 2292
$ws.Range("T5").Value = 12

# Row 6
$ws.Range("D6").Value = 44537
$ws.Range("K6").Value = 'Brooks'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("Q6").Value = '$/caja 20 kilos'
$ws.Range("S6").Value = 1475
$ws.Range("T6").Value = 20

# Row 7
$ws.Range("D7").Value = 44580
$ws.Range("K7").Value = 'Sweet Heart'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 7500
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("S7").Value = 750
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44914
$ws.Range("K8").Value = 'Brooks'
$ws.Range("M8").Value = 700
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 7429
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("S8").Value = 743
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44914
$ws.Range("K9").Value = 'Lapins'
$ws.Range("M9").Value = 550
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7455
$ws.Range("S9").Value = 746

# Row 10
$ws.Range("D10").Value = 44229
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 6500
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6750
$ws.Range("Q10").Value = '$/bandeja 5 kilos'
$ws.Range("R10").Value = 'Provincia de Curicó'
$ws.Range("S10").Value = 1350
$ws.Range("T10").Value = 5

# Row 11
$ws.Range("D11").Value = 44594
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 5000
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 5500
$ws.Range("Q11").Value = '$/bandeja 5 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1100
$ws.Range("T11").Value = 5

# Row 12
$ws.Range("D12").Value = 44908
$ws.Range("K12").Value = 'Rainier'
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15600
$ws.Range("Q12").Value = '$/caja 10 kilos'
$ws.Range("S12").Value = 1560
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("D13").Value = 44571
$ws.Range("K13").Value = 'Brooks'
$ws.Range("M13").Value = 400
$ws.Range("N13").Value = 8500
$ws.Range("O13").Value = 9000
$ws.Range("P13").Value = 8750
$ws.Range("Q13").Value = '$/bandeja 10 kilos'
$ws.Range("S13").Value = 875
$ws.Range("T13").Value = 10

# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("K14").Value = 'Bing'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 39000
$ws.Range("O14").Value = 40000
$ws.Range("P14").Value = 39500
$ws.Range("Q14").Value = '$/caja 20 kilos'
$ws.Range("S14").Value = 1975
$ws.Range("T14").Value = 20

# Row 15
$ws.Range("D15").Value = 44931
$ws.Range("K15").Value = 'Lapins'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6500
$ws.Range("P15").Value = 6250
$ws.Range("S15").Value = 625

# Row 16
$ws.Range("D16").Value = 44931
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3300
$ws.Range("P16").Value = 3150
$ws.Range("Q16").Value = '$/bandeja 5 kilos'
$ws.Range("S16").Value = 630
$ws.Range("T16").Value = 5

# Row 17
$ws.Range("D17").Value = 44568
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("Q17").Value = '$/bandeja 12 kilos'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1292
$ws.Range("T17").Value = 12

# Row 18
$ws.Range("D18").Value = 44210
$ws.Range("K18").Value = 'Rainier'
$ws.Range("N18").Value = 21000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 21500
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("S18").Value = 1194
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44943
$ws.Range("K19").Value = 'Santina'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 600
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 14333
$ws.Range("Q19").Value = '$/caja 15 kilos'
$ws.Range("R19").Value = 'Región del Maule'
$ws.Range("S19").Value = 956
$ws.Range("T19").Value = 15

# Row 20
$ws.Range("D20").Value = 44917
$ws.Range("K20").Value = 'Bing'
$ws.Range("L20").Value = 'Primera'
$ws.Range("N20").Value = 5000
$ws.Range("O20").Value = 6000
$ws.Range("P20").Value = 5625
$ws.Range("S20").Value = 562

# Row 21
$ws.Range("D21").Value = 44917
$ws.Range("K21").Value = 'Santina'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 400
$ws.Range("N21").Value = 5000
$ws.Range("O21").Value = 6000
$ws.Range("P21").Value = 5500
$ws.Range("S21").Value = 550

# Row 22
$ws.Range("D22").Value = 44921
$ws.Range("K22").Value = 'Bing'
$ws.Range("M22").Value = 320
$ws.Range("N22").Value = 7500
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 7781
$ws.Range("Q22").Value = '$/bandeja 10 kilos'
$ws.Range("S22").Value = 778
$ws.Range("T22").Value = 10

# Row 23
$ws.Range("D23").Value = 44901
$ws.Range("M23").Value = 500
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 12500
$ws.Range("Q23").Value = '$/caja 15 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 833
$ws.Range("T23").Value = 15

# Row 24
$ws.Range("D24").Value = 44901
$ws.Range("K24").Value = 'Lapins'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 500
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 12500
$ws.Range("Q24").Value = '$/caja 15 kilos'
$ws.Range("S24").Value = 833
$ws.Range("T24").Value = 15

# Row 25
$ws.Range("D25").Value = 44922
$ws.Range("K25").Value = 'Bing'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 5000
$ws.Range("O25").Value = 6000
$ws.Range("P25").Value = 5500
$ws.Range("Q25").Value = '$/bandeja 10 kilos'
$ws.Range("R25").Value = 'Región del Maule'
$ws.Range("S25").Value = 550
